$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New coalition columns M ("IzquierdaDiv"), N ("Realista"), O ("D2020").
#
# The shared-strings table is built in the literal order new unique strings
# are first written, so we "prime" the table by touching the very first
# occurrence of every brand-new string in the exact sequence required
# before doing the full (idempotent) fill of the three columns.
# ---------------------------------------------------------------------------

$ws.Range("M1").Value  = "IzquierdaDiv"
$ws.Range("M5").Value  = "DC-IC-CIU"
$ws.Range("M10").Value = "Progresistas"
$ws.Range("N1").Value  = "Realista"
$ws.Range("N11").Value = "Chile Digno"
$ws.Range("N5").Value  = "Unidad Constituyente"
$ws.Range("N9").Value  = "Humanista"

# O1 ("D2020") needs a quote-prefixed, mmm-yy-formatted cell (Excel flagged
# the text as looking like a date when it was typed in).
$o1 = $ws.Range("O1")
$o1.Value = "'D2020"
$o1.NumberFormat = "mmm-yy"

$ws.Range("O5").Value = "UC"

# ---------------------------------------------------------------------------
# Full column fill (M, N, O) for every data row.
# ---------------------------------------------------------------------------

$ws.Range("M2").Value  = "CHV"
$ws.Range("N2").Value  = "CHV"
$ws.Range("O2").Value  = "CHV"

$ws.Range("M3").Value  = "CHV"
$ws.Range("N3").Value  = "CHV"
$ws.Range("O3").Value  = "CHV"

$ws.Range("M4").Value  = "CHV"
$ws.Range("N4").Value  = "CHV"
$ws.Range("O4").Value  = "CHV"

$ws.Range("M5").Value  = "DC-IC-CIU"
$ws.Range("N5").Value  = "Unidad Constituyente"
$ws.Range("O5").Value  = "UC"

$ws.Range("M6").Value  = "CHV"
$ws.Range("N6").Value  = "CHV"
$ws.Range("O6").Value  = "CHV"

$ws.Range("M7").Value  = "CHV"
$ws.Range("N7").Value  = "CHV"
$ws.Range("O7").Value  = "CHV"

$ws.Range("M8").Value  = "IND 1"
$ws.Range("N8").Value  = "IND 1"
$ws.Range("O8").Value  = "UC"

$ws.Range("M9").Value  = "IZQ"
$ws.Range("N9").Value  = "Humanista"
$ws.Range("O9").Value  = "Humanista"

$ws.Range("M10").Value = "Progresistas"
$ws.Range("N10").Value = "Unidad Constituyente"
$ws.Range("O10").Value = "UC"

$ws.Range("M11").Value = "IZQ"
$ws.Range("N11").Value = "Chile Digno"
$ws.Range("O11").Value = "PC-FA"

$ws.Range("M12").Value = "IZQ"
$ws.Range("N12").Value = "Chile Digno"
$ws.Range("O12").Value = "PC-FA"

$ws.Range("M13").Value = "PC-FA"
$ws.Range("N13").Value = "FA"
$ws.Range("O13").Value = "UC"

$ws.Range("M14").Value = "PC-FA"
$ws.Range("N14").Value = "FA"
$ws.Range("O14").Value = "PC-FA"

$ws.Range("M15").Value = "PC-FA"
$ws.Range("N15").Value = "Chile Digno"
$ws.Range("O15").Value = "PC-FA"

$ws.Range("M16").Value = "PC-FA"
$ws.Range("N16").Value = "FA"
$ws.Range("O16").Value = "PC-FA"

$ws.Range("M17").Value = "IND"
$ws.Range("N17").Value = "IND"
$ws.Range("O17").Value = "IND"

$ws.Range("M18").Value = "IND 2"
$ws.Range("N18").Value = "Chile Digno"
$ws.Range("O18").Value = "PC-FA"

$ws.Range("M19").Value = "IND 2"
$ws.Range("N19").Value = "Chile Digno"
$ws.Range("O19").Value = "PC-FA"

$ws.Range("M20").Value = "IND 3"
$ws.Range("N20").Value = "IND 3"
$ws.Range("O20").Value = "IND 3"

$ws.Range("M21").Value = "IND 3"
$ws.Range("N21").Value = "IND 3"
$ws.Range("O21").Value = "IND 3"

$ws.Range("M22").Value = "PC-FA"
$ws.Range("N22").Value = "Chile Digno"
$ws.Range("O22").Value = "PC-FA"

$ws.Range("M23").Value = "PC-FA"
$ws.Range("N23").Value = "Chile Digno"
$ws.Range("O23").Value = "PC-FA"

$ws.Range("M24").Value = "DC-IC-CIU"
$ws.Range("N24").Value = "Unidad Constituyente"
$ws.Range("O24").Value = "UC"

$ws.Range("M25").Value = "Progresistas"
$ws.Range("N25").Value = "Unidad Constituyente"
$ws.Range("O25").Value = "UC"

$ws.Range("M26").Value = "Progresistas"
$ws.Range("N26").Value = "Unidad Constituyente"
$ws.Range("O26").Value = "UC"

$ws.Range("M27").Value = "DC-IC-CIU"
$ws.Range("N27").Value = "Unidad Constituyente"
$ws.Range("O27").Value = "UC"

$ws.Range("M28").Value = "Progresistas"
$ws.Range("N28").Value = "Unidad Constituyente"
$ws.Range("O28").Value = "UC"

$ws.Range("M29").Value = "CHV"
$ws.Range("N29").Value = "JAK"
$ws.Range("O29").Value = "JAK"

# ---------------------------------------------------------------------------
# View state: selection moved to P27 (best-effort; matches author's last
# click position before saving).
# ---------------------------------------------------------------------------
$ws.Range("P27").Select()
